$d = $word.ActiveDocument

# Remove the placeholder text "vnpt.SiteAddress" that followed "Địa chỉ: "
$d.Content.Find.Execute("vnpt.SiteAddress", $false, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
